$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("R2").Value = 1 ; $ws.Range("U2").Value = 12
$ws.Range("R6").Value = 1 ; $ws.Range("T6").Value = 1
$ws.Range("R7").Value = 1 ; $ws.Range("U7").Value = 123
$ws.Range("R8").Value = 1 ; $ws.Range("U8").Value = 12
$ws.Range("R9").Value = 1 ; $ws.Range("T9").Value = 1
$ws.Range("S10").Value = 1 ; $ws.Range("U10").Value = 13
$ws.Range("R11").Value = 1 ; $ws.Range("S11").Value = 1 ; $ws.Range("U11").Value = 123
$ws.Range("R12").Value = 1 ; $ws.Range("U12").Value = 123
$ws.Range("R14").Value = 1 ; $ws.Range("U14").Value = 12
$ws.Range("R15").Value = 1 ; $ws.Range("U15").Value = 12
$ws.Range("R16").Value = 1 ; $ws.Range("U16").Value = 12
$ws.Range("R17").Value = 1 ; $ws.Range("U17").Value = 123
$ws.Range("R18").Value = 1 ; $ws.Range("U18").Value = 123
$ws.Range("Q19").Value = 1 ; $ws.Range("S19").Value = 1 ; $ws.Range("T19").Value = 1 ; $ws.Range("U19").Value = 13
$ws.Range("Q21").Value = 1 ; $ws.Range("R21").Value = 1 ; $ws.Range("T21").Value = 1 ; $ws.Range("U21").Value = 12
$ws.Range("R22").Value = 1 ; $ws.Range("U22").Value = 123
$ws.Range("R23").Value = 1 ; $ws.Range("T23").Value = 1
$ws.Range("R24").Value = 1 ; $ws.Range("T24").Value = 1
$ws.Range("Q25").Value = 1 ; $ws.Range("S25").Value = 1 ; $ws.Range("T25").Value = 1 ; $ws.Range("U25").Value = 13
$ws.Range("R26").Value = 1 ; $ws.Range("U26").Value = 123
$ws.Range("Q27").Value = 1 ; $ws.Range("R27").Value = 1 ; $ws.Range("S27").Value = 1 ; $ws.Range("T27").Value = 1 ; $ws.Range("U27").Value = 123
$ws.Range("Q28").Value = 1 ; $ws.Range("R28").Value = 1 ; $ws.Range("T28").Value = 1 ; $ws.Range("U28").Value = 12
$ws.Range("R29").Value = 1 ; $ws.Range("T29").Value = 1
$ws.Range("R30").Value = 1 ; $ws.Range("U30").Value = 123
$ws.Range("R32").Value = 1 ; $ws.Range("T32").Value = 1
$ws.Range("Q33").Value = 1 ; $ws.Range("S33").Value = 1 ; $ws.Range("T33").Value = 1 ; $ws.Range("U33").Value = 13
$ws.Range("S34").Value = 1 ; $ws.Range("T34").Value = 1
$ws.Range("R37").Value = 1 ; $ws.Range("U37").Value = 12
$ws.Range("Q41").Value = 1 ; $ws.Range("R41").Value = 1 ; $ws.Range("T41").Value = 1 ; $ws.Range("U41").Value = 12
$ws.Range("R42").Value = 1 ; $ws.Range("U42").Value = 123
$ws.Range("S47").Value = 1 ; $ws.Range("U47").Value = 13
$ws.Range("Q48").Value = 1 ; $ws.Range("U48").Value = 13
$ws.Range("Q49").Value = 1 ; $ws.Range("S49").Value = 1 ; $ws.Range("T49").Value = 1 ; $ws.Range("U49").Value = 13
$ws.Range("R50").Value = 1 ; $ws.Range("U50").Value = 12
$ws.Range("R51").Value = 1 ; $ws.Range("T51").Value = 1
$ws.Range("S53").Value = 1 ; $ws.Range("T53").Value = 1
$ws.Range("S54").Value = 1 ; $ws.Range("T54").Value = 1
$ws.Range("R56").Value = 1 ; $ws.Range("U56").Value = 12
$ws.Range("R57").Value = 1 ; $ws.Range("U57").Value = 12
$ws.Range("Q59").Value = 1 ; $ws.Range("R59").Value = 1 ; $ws.Range("T59").Value = 1 ; $ws.Range("U59").Value = 12
$ws.Range("R60").Value = 1 ; $ws.Range("U60").Value = 12
$ws.Range("R61").Value = 1 ; $ws.Range("U61").Value = 123
$ws.Range("R62").Value = 1 ; $ws.Range("T62").Value = 1
$ws.Range("R63").Value = 1 ; $ws.Range("T63").Value = 1
$ws.Range("R65").Value = 1 ; $ws.Range("T65").Value = 1
$ws.Range("R66").Value = 1 ; $ws.Range("U66").Value = 12
$ws.Range("Q67").Value = 1 ; $ws.Range("R67").Value = 1 ; $ws.Range("T67").Value = 1 ; $ws.Range("U67").Value = 12
$ws.Range("R69").Value = 1 ; $ws.Range("S69").Value = 1 ; $ws.Range("T69").Value = 1 ; $ws.Range("U69").Value = 23
$ws.Range("R70").Value = 1 ; $ws.Range("U70").Value = 123
$ws.Range("R72").Value = 1 ; $ws.Range("T72").Value = 1
$ws.Range("Q73").Value = 1 ; $ws.Range("R73").Value = 1 ; $ws.Range("T73").Value = 1 ; $ws.Range("U73").Value = 12
$ws.Range("Q74").Value = 1 ; $ws.Range("R74").Value = 1 ; $ws.Range("S74").Value = 1 ; $ws.Range("T74").Value = 1 ; $ws.Range("U74").Value = 123
$ws.Range("S77").Value = 1 ; $ws.Range("T77").Value = 1
$ws.Range("R79").Value = 1 ; $ws.Range("U79").Value = 12
$ws.Range("R80").Value = 1 ; $ws.Range("U80").Value = 12
$ws.Range("Q81").Value = 1 ; $ws.Range("S81").Value = 1 ; $ws.Range("T81").Value = 1 ; $ws.Range("U81").Value = 13
$ws.Range("Q82").Value = 1 ; $ws.Range("R82").Value = 1 ; $ws.Range("T82").Value = 1 ; $ws.Range("U82").Value = 12
$ws.Range("R83").Value = 1 ; $ws.Range("T83").Value = 1
$ws.Range("R84").Value = 1 ; $ws.Range("S84").Value = 1 ; $ws.Range("U84").Value = 123
$ws.Range("R85").Value = 1 ; $ws.Range("U85").Value = 12
$ws.Range("R86").Value = 1 ; $ws.Range("T86").Value = 1
$ws.Range("R87").Value = 1 ; $ws.Range("T87").Value = 1
$ws.Range("R88").Value = 1 ; $ws.Range("U88").Value = 12
$ws.Range("R89").Value = 1 ; $ws.Range("S89").Value = 1 ; $ws.Range("U89").Value = 123
$ws.Range("R90").Value = 1 ; $ws.Range("S90").Value = 1 ; $ws.Range("U90").Value = 123
$ws.Range("R91").Value = 1 ; $ws.Range("T91").Value = 1
$ws.Range("R92").Value = 1 ; $ws.Range("S92").Value = 1 ; $ws.Range("U92").Value = 123
$ws.Range("R93").Value = 1 ; $ws.Range("U93").Value = 12
$ws.Range("R94").Value = 1 ; $ws.Range("U94").Value = 12
$ws.Range("R95").Value = 1 ; $ws.Range("T95").Value = 1
$ws.Range("Q98").Value = 1 ; $ws.Range("R98").Value = 1 ; $ws.Range("T98").Value = 1 ; $ws.Range("U98").Value = 12
$ws.Range("S100").Value = 1 ; $ws.Range("T100").Value = 1
$ws.Range("R101").Value = 1 ; $ws.Range("T101").Value = 1
$ws.Range("R102").Value = 1 ; $ws.Range("S102").Value = 1 ; $ws.Range("T102").Value = 1 ; $ws.Range("U102").Value = 23
$ws.Range("Q104").Value = 1 ; $ws.Range("S104").Value = 1 ; $ws.Range("T104").Value = 1 ; $ws.Range("U104").Value = 13
$ws.Range("R107").Value = 1 ; $ws.Range("U107").Value = 123
$ws.Range("R108").Value = 1 ; $ws.Range("U108").Value = 123
$ws.Range("Q109").Value = 1 ; $ws.Range("R109").Value = 1 ; $ws.Range("T109").Value = 1 ; $ws.Range("U109").Value = 12
$ws.Range("S110").Value = 1 ; $ws.Range("T110").Value = 1
$ws.Range("Q112").Value = 1 ; $ws.Range("R112").Value = 1 ; $ws.Range("S112").Value = 1 ; $ws.Range("T112").Value = 1 ; $ws.Range("U112").Value = 123
$ws.Range("R113").Value = 1 ; $ws.Range("T113").Value = 1
